# Data was modified in the IPL folder: two new columns (ownTeam, oppTeam)
# were inserted after "result" and before "batsman", and the match log was
# expanded from 5 rows to a full 13-row history (header unaffected aside
# from the two new columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('venue', 'date', 'result', 'ownTeam', 'oppTeam', 'batsman', 'totalRuns', 'totalBalls', 'total4s', 'total6s', 'sr'),
    @(' Sharjah', ' September 27 2020', 'Royals won by 4 wickets (with 3 balls remaining)', 'Kings XI Punjab', 'Rajasthan Royals', 'Nicholas Pooran ', '25', '8', '1', '3', '312.50'),
    @(' Abu Dhabi', ' October 30 2020', 'Royals won by 7 wickets (with 15 balls remaining)', 'Kings XI Punjab', 'Rajasthan Royals', 'Nicholas Pooran ', '22', '10', '0', '3', '220.00'),
    @(' Dubai (DSC)', ' October 24 2020', 'Kings XI won by 12 runs', 'Kings XI Punjab', 'Sunrisers Hyderabad', 'Nicholas Pooran ', '32', '28', '2', '0', '114.28'),
    @(' Dubai (DSC)', ' September 24 2020', 'Kings XI won by 97 runs', 'Kings XI Punjab', 'Royal Challengers Bangalore', 'Nicholas Pooran ', '17', '18', '1', '0', '94.44'),
    @(' Sharjah', ' October 15 2020', 'Kings XI won by 8 wickets', 'Kings XI Punjab', 'Royal Challengers Bangalore', 'Nicholas Pooran ', '6', '1', '0', '1', '600.00'),
    @(' Abu Dhabi', ' November 01 2020', 'Super Kings won by 9 wickets (with 7 balls remaining)', 'Kings XI Punjab', 'Chennai Super Kings', 'Nicholas Pooran ', '2', '6', '0', '0', '33.33'),
    @(' Abu Dhabi', ' October 01 2020', 'Mumbai won by 48 runs', 'Kings XI Punjab', 'Mumbai Indians', 'Nicholas Pooran ', '44', '27', '3', '2', '162.96'),
    @(' Abu Dhabi', ' October 10 2020', 'KKR won by 2 runs', 'Kings XI Punjab', 'Kolkata Knight Riders', 'Nicholas Pooran ', '16', '10', '2', '1', '160.00'),
    @(' Dubai (DSC)', ' October 20 2020', 'Kings XI won by 5 wickets (with 6 balls remaining)', 'Kings XI Punjab', 'Delhi Capitals', 'Nicholas Pooran ', '53', '28', '6', '3', '189.28'),
    @(' Dubai (DSC)', ' September 20 2020', 'Match tied (Capitals won the one-over eliminator)', 'Kings XI Punjab', 'Delhi Capitals', 'Nicholas Pooran ', '0', '3', '0', '0', '0.00'),
    @(' Sharjah', ' October 26 2020', 'Kings XI won by 8 wickets (with 7 balls remaining)', 'Kings XI Punjab', 'Kolkata Knight Riders', 'Nicholas Pooran ', '2', '3', '0', '0', '66.66'),
    @(' Dubai (DSC)', ' October 18 2020', 'Match tied (Kings XI won the one-over eliminator)', 'Kings XI Punjab', 'Mumbai Indians', 'Nicholas Pooran ', '24', '12', '2', '2', '200.00'),
    @(' Dubai (DSC)', ' October 08 2020', 'Sunrisers won by 69 runs', 'Kings XI Punjab', 'Sunrisers Hyderabad', 'Nicholas Pooran ', '77', '37', '5', '7', '208.10')
)

$rowCount = $data.Count
$colCount = $data[0].Count

# Columns G:K (totalRuns, totalBalls, total4s, total6s, sr) hold numeric-looking
# text (e.g. "312.50", "220.00") that must stay text, not become numbers.
$ws.Range("G2:K" + $rowCount).NumberFormat = "@"

for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}
